$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 - shifts existing rows 15-17 down to 16-18
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the latest week's data
# (same market/category/etc. as the surrounding rows, new date + prices)
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44943
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100114007
$ws.Cells.Item(15, 7).Value = "Jengibre"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 17000
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 17000
$ws.Cells.Item(15, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 1308
$ws.Cells.Item(15, 17).Value = 13
$ws.Cells.Item(15, 18).Value = "Hortaliza"
